$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (G1:I1), matching the bold/centered/bordered style of the existing header cells
$headers = @{
    "G1" = "HISTOGRAM_VARIANCE"
    "H1" = "NPCR"
    "I1" = "UACI"
}
foreach ($addr in $headers.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $headers[$addr]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Row data: image name (A) plus metric columns B-I (rows re-sorted alphabetically by filename,
# and three new metric columns added: HISTOGRAM_VARIANCE, NPCR, UACI)
$data = @(
    @(2, "6.1.01.tiff", 8.156778997588813, 105.0476531982422, 0.007559322763799163, 7.993019045411336, 0.0002945513075038, 641.116124567474, 99.66888427734375, 33.41211655560662),
    @(3, "Aerial.tiff", 9.23117716342462, 105.6394500732422, 0.007969840369582818, 7.995576246553931, -0.004108210209782912, 404.7654286812764, 99.56207275390625, 33.62072514552696),
    @(4, "Aerial512.tiff", 8.173851793941839, 105.5501480102539, 0.007912867873204462, 7.993844373922263, -0.001089993672934644, 9083.962660515186, 99.50180053710938, 33.37418799306832),
    @(5, "Airplane.tiff", 7.700057863821832, 105.4205017089844, 0.008077793673996381, 7.98969225199809, -0.001681104114651707, 947.2059054209919, 99.578857421875, 33.65193684895834),
    @(6, "Airplane512.tiff", 8.838424511909482, 105.4510612487793, 0.009745887063940707, 7.931540414743385, -0.001784337156529359, 100261.8439061899, 99.57313537597656, 34.3744150797526),
    @(7, "APC512.tiff", 9.772740090814745, 105.5988731384277, 0.0103780066829179, 7.986017729319951, 0.0008936406818229925, 20516.7632141484, 99.38392639160156, 33.53911007151884),
    @(8, "boat.512.tiff", 9.275575313728732, 105.4716911315918, 0.01000876096399781, 7.995354611728661, 0.003687502162863957, 6775.368242983468, 99.50141906738281, 33.28296586578968),
    @(9, "CarAPCs1_512.tiff", 10.26772769313258, 105.6706886291504, 0.01016499160758506, 7.955389150778685, -0.0003471934357054807, 65789.20870434448, 99.32403564453125, 33.66560244092754),
    @(10, "CarAPCs_512.tiff", 10.12962306499437, 105.5632247924805, 0.01031677715744493, 7.983693310550926, 0.002046856030146311, 24045.89004229143, 99.39422607421875, 33.43713199391085),
    @(11, "Chemical_plant.tiff", 9.212167315570367, 105.563720703125, 0.005446752314781257, 7.996077428765373, -0.007380787365136933, 357.6086735870819, 99.56817626953125, 33.14146752450981),
    @(12, "Clock.tiff", 7.223503761205917, 106.2551879882812, 0.007603752979439478, 7.988889079488572, -0.0007939449363241238, 1028.393633217993, 99.43084716796875, 33.2638370289522),
    @(13, "Couple.tiff", 9.593233710095811, 105.5767555236816, 0.009287104179924878, 7.996393633895227, -0.001580719728598491, 5300.115340253748, 99.560546875, 33.52138294893152),
    @(14, "gray21.512.tiff", 7.581761094389405, 105.4107704162598, 0.008125846554992003, 7.967006975055557, 0.00269277671513107, 48331.43781622453, 99.55024719238281, 33.31559573902803),
    @(15, "Moon_256.tiff", 10.11697339965389, 105.7922058105469, 0.01149913499084049, 7.995315195510237, 0.002541187748652892, 426.080061514802, 99.462890625, 33.62401625689338),
    @(16, "Resolution_chart.tiff", 5.075455745271743, 105.5915222167969, 0.007517480493097176, 7.731313578838155, 0.001435207900537723, 27466.42798923491, 99.9542236328125, 46.07293820848652),
    @(17, "ruler.512.tiff", 4.809952977881065, 107.4700088500977, 0.03242234337353897, 7.367217917872677, 0.04059868488305154, 1027000.308189158, 100, 45.09413924871706),
    @(18, "StreamandBridge512.tiff", 8.765354531848962, 105.4198379516602, 0.009261770987249804, 7.980706517660061, 0.000832837947973298, 28343.48536716648, 99.29237365722656, 33.41863893995098),
    @(19, "Tank1_512.tiff", 10.03672480097016, 105.7082977294922, 0.009367194069846192, 7.986213130618101, -0.0005680859479606264, 20191.6186082276, 99.45182800292969, 33.24813693177466),
    @(20, "Tank2_512.tiff", 9.77020087941524, 105.2868003845215, 0.009424441825384112, 7.989410828144559, -0.001267870697635909, 15572.62665128797, 99.066162109375, 33.19874782188266),
    @(21, "Tank512.tiff", 10.09596674981244, 105.3793525695801, 0.009606588336180808, 7.976850258422874, -0.001493196028722689, 33868.32166089965, 99.03945922851562, 33.43890620212929),
    @(22, "Truck512.tiff", 9.940937438829923, 105.6881713867188, 0.01026499710230641, 7.98644140774225, 0.00130869452956219, 19691.42945021146, 99.43313598632812, 33.45802007936964),
    @(23, "TruckAPCs1_512.tiff", 9.137623429600824, 105.5090560913086, 0.009540003759556928, 7.992894890550668, 0.001323941815422989, 10410.02100730488, 99.530029296875, 33.36989010081572),
    @(24, "TruckAPCs512.tiff", 9.570980281418064, 105.5757217407227, 0.01053406707997141, 7.991691195245294, 0.0007599210695027543, 12183.98585159554, 99.49302673339844, 33.52707956351486),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}
